$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.386581420898438
$ws.Range("B1").Value = 1.541714668273926
$ws.Range("C1").Value = 5.142643928527832
$ws.Range("D1").Value = 2.737042188644409
$ws.Range("E1").Value = 0.9275162220001221
